# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 15.12 = 61755.73 pesos"), "1000 Bs = 15.15 = 61920.48 pesos"
$text = $text -replace [regex]::Escape("61755.73 pesos = 15.0 = 968.91 Bs"), "61920.48 pesos = 15.07 = 981.09 Bs"
$cell.Value = $text

# --- tasas: update the rate table values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 65.999
$wsTasas.Range("O10").Value = 4086.69
$wsTasas.Range("N12").Value = 4110
$wsTasas.Range("O12").Value = 65.12
